$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking values (e.g. "109.14") are not
# auto-converted to numbers by Excel, matching the existing inline-string format.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "47.287.84"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "2.501.84"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("D5").Value = "323.68"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "109.14"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "39.18"
$ws.Range("E10").Value = "  +8.55%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "18.40"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "2.892.24"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "2.502.55"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "47.217.26"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").Value = "12.87"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "6.64"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  +12.60%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "247.95"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("D26").Value = "26.08"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +5.01%  "
$ws.Range("D29").Value = "10.08"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("D30").Value = "35.24"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "0.139"
$ws.Range("E31").Value = "  +7.24%  "
$ws.Range("D32").Value = "49.83"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "20.05"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "0.0790"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("D38").Value = "4.71"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "122.33"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "21.24"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "1.992.02"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").Value = "56.54"
$ws.Range("E51").Value = "  +3.56%  "
